# Updated with new versions of login, and search
#
# This script reproduces (as closely as the COM surface allows) the authored
# change: a brand-new worksheet named "Shay" was added after "default",
# populated with a (partial) copy of the data/headers that live on
# "default" (rows 1-4, columns A:N), plus blank-but-styled filler rows
# down through row 20, some column widths / zoom tweaks, and a new
# selection on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1) Add the new sheet "Shay" directly after the existing "default" tab.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Shay"

# ------------------------------------------------------------------
# 2) Populate "Shay" with the same header/data block as "default"
#    (A1:N4 -> header row + 3 data rows), carrying over the exact
#    same cell styles (fills/fonts/number formats) via a Range copy
#    instead of re-creating styles by hand.
# ------------------------------------------------------------------
$ws1.Range("A1:N4").Copy($ws2.Range("A1"))

# Blank styled filler rows 5:6 across A:N (mirrors default's blank
# body rows, including the grey separator in column F).
$ws1.Range("A5:N5").Copy($ws2.Range("A5:N6"))

# Blank styled filler rows 7:20 across A:L only (no M/N on these rows).
$ws1.Range("A5:L5").Copy($ws2.Range("A7:L20"))

# ------------------------------------------------------------------
# 3) Cosmetic sheet-level tweaks on "Shay".
# ------------------------------------------------------------------
$ws2.Range("A1:L20").RowHeight = 21
$ws2.Range("A1:L12").ColumnWidth = 19.877604166666668

$ws2.Activate()
$excel.ActiveWindow.Zoom = 85
$ws2.Range("E8").Select()

# ------------------------------------------------------------------
# 4) "default" loses the tab-selected flag (handled automatically by
#    activating "Shay" above) and its selection becomes A1:E2.
# ------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A1:E2").Select()

# Leave "Shay" as the active/selected sheet, matching activeTab="1".
$ws2.Activate()
